$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B-E)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values (columns B-E)
$ws.Range("B2").Value = 44.519599243087939
$ws.Range("C2").Value = -8.9276162598221021
$ws.Range("D2").Value = -3.3074225400202302
$ws.Range("E2").Value = 6.1040602419814824

# Row 3 values - C3 is cleared, D3 gets a new value
$ws.Range("B3").Value = 44.130543730790535
$ws.Range("C3").ClearContents() | Out-Null
$ws.Range("D3").Value = -20.962144060874568
$ws.Range("E3").Value = 22.123273101918144

# Update the selection to match the new used range
$ws.Range("B1:E3").Select() | Out-Null
